# Fix typo in Cheat sheet with regex
#
# The "except the string 'Species'" row of the regex-examples table on
# slide 1 had a doubled leading single-quote in its code sample
# ( ''^(?!Species$).*'  ->  '^(?!Species$).*' ). Locate that table cell
# through the normal Shape/Table/Cell object model and correct the text
# while leaving every other run property (font, size, etc.) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Table 84") {
        $targetShape = $shp
        break
    }
}

$tbl = $targetShape.Table

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $cellShape = $tbl.Rows.Item($r).Cells.Item(1).Shape
    $tr = $cellShape.TextFrame.TextRange
    if ($tr.Text -eq "''^(?!Species`$).*'") {
        $tr.Text = "'^(?!Species`$).*'"
    }
}
